# fix: typos on menu
# - "rødløk" -> "raudlauk" in the sandwich description
# - move the saved selection/active cell to C10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C4")
$cell.Value = "aioli, salat, agurk, paprika, raudlauk og tomat. 1,2,3,10,12"

$ws.Range("C10").Select()
